$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = 0.07581796488389111
$ws.Range("Y2").Value = 0.2363100966135811
$ws.Range("AB2").Value = 0.07578890859571075
$ws.Range("AC2").Value = -0.9549658756721473

$ws.Range("X3").Value = 0.07583597537525198
$ws.Range("Y3").Value = 0.1713550358607031
$ws.Range("AB3").Value = 0.07580429214353542
$ws.Range("AC3").Value = -0.2175726596516978

$ws.Range("X4").Value = 0.07579995439253025
$ws.Range("Y4").Value = 0.301265157366459
$ws.Range("AB4").Value = 0.07577352504788608
$ws.Range("AC4").Value = -1.692359091692597
